# Update NATMI LR-pair sheet (L1cam-Itgav) with new TPM-based expression
# values. The "Ligand/Receptor average & total expression value" columns
# (G/H and M/N) are driven by per-cluster TPM numbers; every other edited
# column (I/J, O/P specificity and Q/R/S/T edge-weight columns) is derived
# from those base numbers, so we recompute the whole dependent chain here
# instead of poking in literal numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new (post "new tpm") per-cluster ligand average expression values ----
# (column G for rows whose Sending cluster == this cluster)
$newG = @{}
$newG["ECs"]               = 5.375839
$newG["FAPs"]               = 0.1628146666666667
$newG["Inflammatory-Mac"]   = 9.994147
$newG["MuSCs"]              = 0.7761303333333333
$newG["Resolving-Mac"]      = 6.522593333333333

# ---- new per-cluster receptor average expression values ----
# (column M for rows whose Target cluster == this cluster)
$newM = @{}
$newM["ECs"]               = 22.495411
$newM["FAPs"]               = 82.64333833333332
$newM["Inflammatory-Mac"]   = 79.32606499999999
$newM["MuSCs"]              = 14.467164
$newM["Resolving-Mac"]      = 72.34725666666667

# total-cell multiplier used to turn an average expression value into a
# total expression value (Ligand/Receptor total expression value columns)
$cellCount = 3

$clusters = @("ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Resolving-Mac")

# Ligand/Receptor derived-specificity denominators: sum of the average
# expression value across all clusters
$sumG = 0.0
foreach ($k in $clusters) { $sumG = $sumG + $newG[$k] }

$sumM = 0.0
foreach ($k in $clusters) { $sumM = $sumM + $newM[$k] }

# first pass: write G,H,I,J,M,N,O,P and keep the edge-weight numbers
# (Q = G*M, R = H*N) around so we can normalise them (S,T) in a second
# pass once their grand totals are known
$rows = @()
$sumQ = 0.0
$sumR = 0.0

for ($si = 0; $si -lt $clusters.Length; $si++) {
    $sendCluster = $clusters[$si]
    $G = $newG[$sendCluster]
    $H = $G * $cellCount
    $I = $G / $sumG
    $J = $I

    for ($ti = 0; $ti -lt $clusters.Length; $ti++) {
        $targetCluster = $clusters[$ti]
        $M = $newM[$targetCluster]
        $N = $M * $cellCount
        $O = $M / $sumM
        $P = $O

        $Q = $G * $M
        $R = $H * $N

        $row = 2 + ($si * $clusters.Length) + $ti

        $ws.Range("G$row").Value = $G
        $ws.Range("H$row").Value = $H
        $ws.Range("I$row").Value = $I
        $ws.Range("J$row").Value = $J

        $ws.Range("M$row").Value = $M
        $ws.Range("N$row").Value = $N
        $ws.Range("O$row").Value = $O
        $ws.Range("P$row").Value = $P

        $ws.Range("Q$row").Value = $Q
        $ws.Range("R$row").Value = $R

        $sumQ = $sumQ + $Q
        $sumR = $sumR + $R

        $rows += , @($row, $Q, $R)
    }
}

# second pass: edge-weight derived specificity = value / grand total
foreach ($entry in $rows) {
    $row = $entry[0]
    $Q = $entry[1]
    $R = $entry[2]

    $S = $Q / $sumQ
    $T = $R / $sumR

    $ws.Range("S$row").Value = $S
    $ws.Range("T$row").Value = $T
}

Write-Host "Updated rows 2-26 (G,H,I,J,M,N,O,P,Q,R,S,T) with new TPM values"
